$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1362
$ws.Range("C3").Value = 682
$ws.Range("C4").Value = 540
$ws.Range("C5").Value = 189
$ws.Range("C6").Value = 189
$ws.Range("C7").Value = 170
$ws.Range("C8").Value = 170
$ws.Range("C9").Value = 134
$ws.Range("C10").Value = 134
$ws.Range("C11").Value = 132
